$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2;  Date = "2025/11/22"; Val = "4.42" },
    @{ Row = 8;  Date = "2025/11/22"; Val = "7.41" },
    @{ Row = 14; Date = "2025/11/22"; Val = "2.76" },
    @{ Row = 20; Date = "2025/11/22"; Val = "12.20" },
    @{ Row = 26; Date = "2025/11/22"; Val = "9.71" },
    @{ Row = 32; Date = "2025/11/22"; Val = "25.34" },
    @{ Row = 38; Date = "2025/11/22"; Val = $null },
    @{ Row = 44; Date = "2025/11/22"; Val = "9.94" },
    @{ Row = 50; Date = "2025/11/22"; Val = "11.04" },
    @{ Row = 56; Date = "2025/11/22"; Val = "29.34" },
    @{ Row = 62; Date = "2025/11/22"; Val = "10.75" },
    @{ Row = 68; Date = "2025/11/22"; Val = "11.63" },
    @{ Row = 74; Date = "2025/11/22"; Val = "15.13" }
)

foreach ($item in $rows) {
    $r = $item.Row

    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $item.Date

    if ($item.Val -ne $null) {
        $cellB = $ws.Cells.Item($r, 2)
        $cellB.NumberFormat = "@"
        $cellB.Value = $item.Val
    }
}
